# Updated cryptos list - apply latest price/volume figures as plain text cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay text (matches source data formatting such as
    # "30.648.16" or "1.001") instead of being auto-coerced to a number,
    # then restore the default "Normal" style so no stray number-format
    # style gets attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.648.16"
Set-TextValue $ws.Range("E2") "  +1.36%  "
Set-TextValue $ws.Range("D3") "1.890.59"
Set-TextValue $ws.Range("E3") "  +1.74%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("E4") "  +0.20%  "
Set-TextValue $ws.Range("D5") "238.27"
Set-TextValue $ws.Range("E5") "  +1.06%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.17%  "
Set-TextValue $ws.Range("D7") "0.4830"
Set-TextValue $ws.Range("E7") "  +1.02%  "
Set-TextValue $ws.Range("D8") "0.2875"
Set-TextValue $ws.Range("E8") "  +2.42%  "
Set-TextValue $ws.Range("D9") "0.06550"
Set-TextValue $ws.Range("E9") "  +1.69%  "
Set-TextValue $ws.Range("D10") "1.867.23"
Set-TextValue $ws.Range("E10") "  +0.48%  "
Set-TextValue $ws.Range("D11") "16.78"
Set-TextValue $ws.Range("E11") "  +3.34%  "
Set-TextValue $ws.Range("D12") "0.07453"
Set-TextValue $ws.Range("E12") "  +0.84%  "
Set-TextValue $ws.Range("D13") "5.103"
Set-TextValue $ws.Range("E13") "  +0.15%  "
Set-TextValue $ws.Range("D14") "87.85"
Set-TextValue $ws.Range("E14") "  +0.85%  "
Set-TextValue $ws.Range("D15") "0.6659"
Set-TextValue $ws.Range("E15") "  +3.29%  "
Set-TextValue $ws.Range("D16") "30.613.97"
Set-TextValue $ws.Range("E16") "  +1.48%  "
Set-TextValue $ws.Range("D17") "13.23"
Set-TextValue $ws.Range("E17") "  +0.52%  "
Set-TextValue $ws.Range("D18") "1.001"
Set-TextValue $ws.Range("E18") "  +0.14%  "
Set-TextValue $ws.Range("D19") "0.000007583"
Set-TextValue $ws.Range("E19") "  +0.30%  "
Set-TextValue $ws.Range("D20") "231.87"
Set-TextValue $ws.Range("E20") "  +3.45%  "
Set-TextValue $ws.Range("D21") "1.002"
Set-TextValue $ws.Range("E21") "  +0.19%  "
Set-TextValue $ws.Range("D22") "5.273"
Set-TextValue $ws.Range("E22") "  +0.00%  "
Set-TextValue $ws.Range("D23") "2.038.37"
Set-TextValue $ws.Range("E23") "  -2.82%  "
Set-TextValue $ws.Range("D24") "6.187"
Set-TextValue $ws.Range("E24") "  +1.71%  "
Set-TextValue $ws.Range("D25") "9.383"
Set-TextValue $ws.Range("E25") "  +1.95%  "
Set-TextValue $ws.Range("D26") "168.61"
Set-TextValue $ws.Range("E26") "  +3.02%  "
Set-TextValue $ws.Range("D27") "18.71"
Set-TextValue $ws.Range("E27") "  +1.18%  "
Set-TextValue $ws.Range("D28") "1.956"
Set-TextValue $ws.Range("E28") "  +1.61%  "
Set-TextValue $ws.Range("E29") "  +11.39%  "
Set-TextValue $ws.Range("E30") "  -2.35%  "
Set-TextValue $ws.Range("D31") "4.329"
Set-TextValue $ws.Range("E31") "  +2.24%  "
Set-TextValue $ws.Range("D32") "4.022"
Set-TextValue $ws.Range("E32") "  +1.94%  "
Set-TextValue $ws.Range("D33") "0.05066"
Set-TextValue $ws.Range("E33") "  +1.83%  "
Set-TextValue $ws.Range("D34") "1.209"
Set-TextValue $ws.Range("E34") "  +5.54%  "
Set-TextValue $ws.Range("D35") "0.7505"
Set-TextValue $ws.Range("E35") "  +3.80%  "
Set-TextValue $ws.Range("D36") "0.9989"
Set-TextValue $ws.Range("E36") "  +0.16%  "
Set-TextValue $ws.Range("D37") "2.713"
Set-TextValue $ws.Range("E37") "  +1.03%  "
Set-TextValue $ws.Range("E38") "  +3.41%  "
Set-TextValue $ws.Range("D39") "2.653"
Set-TextValue $ws.Range("E39") "  +2.26%  "
Set-TextValue $ws.Range("D40") "0.9206"
Set-TextValue $ws.Range("E40") "  +2.36%  "
Set-TextValue $ws.Range("D41") "2.063"
Set-TextValue $ws.Range("E41") "  +1.41%  "
Set-TextValue $ws.Range("D42") "106.72"
Set-TextValue $ws.Range("E42") "  +0.90%  "
Set-TextValue $ws.Range("D43") "0.4284"
Set-TextValue $ws.Range("E43") "  +0.87%  "
Set-TextValue $ws.Range("E44") "  +0.30%  "
Set-TextValue $ws.Range("D45") "5.629"
Set-TextValue $ws.Range("E45") "  -4.02%  "
Set-TextValue $ws.Range("D46") "7.418"
Set-TextValue $ws.Range("E46") "  +1.99%  "
Set-TextValue $ws.Range("D47") "64.17"
Set-TextValue $ws.Range("E47") "  +0.63%  "
Set-TextValue $ws.Range("D48") "0.1277"
Set-TextValue $ws.Range("E48") "  -2.08%  "
Set-TextValue $ws.Range("D49") "1.487"
Set-TextValue $ws.Range("E49") "  -0.45%  "
Set-TextValue $ws.Range("D50") "9.017"
Set-TextValue $ws.Range("E50") "  +4.00%  "
Set-TextValue $ws.Range("D51") "34.00"
Set-TextValue $ws.Range("E51") "  +1.08%  "
